$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 107; this shifts all rows 107-190 down to 108-191
$ws.Rows(107).Insert()

# Populate the newly inserted row 107 with the new record's data
$ws.Range("A107").Value = 5
$ws.Range("B107").Value = "Macroferia Regional de Talca"
$ws.Range("C107").Value = "Maule"
$ws.Range("D107").Value = 44977
$ws.Range("E107").Value = 7
$ws.Range("F107").Value = 100112030
$ws.Range("G107").Value = "Poroto granado"
$ws.Range("H107").Value = "Sin especificar"
$ws.Range("I107").Value = "Primera"
$ws.Range("J107").Value = 300
$ws.Range("K107").Value = 25000
$ws.Range("L107").Value = 25000
$ws.Range("M107").Value = 25000
$ws.Range("N107").Value = "$/saco 25 kilos"
$ws.Range("O107").Value = "Región del Maule"
$ws.Range("P107").Value = 1000
$ws.Range("Q107").Value = 25
$ws.Range("R107").Value = "Hortaliza"
